# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# sheets, which carry duplicate data.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8004
    3  = 7592
    6  = 37
    8  = 125
    9  = 105
    10 = 147
    11 = 222
    12 = 686
    13 = 109
    14 = 1180
    15 = 56
    19 = 97
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
